$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column R (col 18) with year 2021 header and data, mirroring
# existing columns' formatting where applicable.

# Header: R3 = 2021 (numeric), same style as Q3
$ws.Cells.Item(3, 18).Value = 2021
$ws.Cells.Item(3, 18).Style = $ws.Cells.Item(3, 17).Style

# Row 4: numeric value, style copied from N4/O4/P4 (numFmt 164 style)
$ws.Cells.Item(4, 18).Value = 0.00029886145739191973
$ws.Cells.Item(4, 18).Style = $ws.Cells.Item(4, 16).Style

# Rows 5-11: "-" text, same style as Q column in those rows
$ws.Cells.Item(5, 18).Value = "-"
$ws.Cells.Item(5, 18).Style = $ws.Cells.Item(5, 17).Style

$ws.Cells.Item(6, 18).Value = "-"
$ws.Cells.Item(6, 18).Style = $ws.Cells.Item(6, 17).Style

$ws.Cells.Item(7, 18).Value = "-"
$ws.Cells.Item(7, 18).Style = $ws.Cells.Item(7, 17).Style

$ws.Cells.Item(8, 18).Value = "-"
$ws.Cells.Item(8, 18).Style = $ws.Cells.Item(8, 17).Style

$ws.Cells.Item(9, 18).Value = "-"
$ws.Cells.Item(9, 18).Style = $ws.Cells.Item(9, 17).Style

$ws.Cells.Item(10, 18).Value = "-"
$ws.Cells.Item(10, 18).Style = $ws.Cells.Item(10, 17).Style

$ws.Cells.Item(11, 18).Value = "-"
$ws.Cells.Item(11, 18).Style = $ws.Cells.Item(11, 17).Style

# Row 12: numeric value with a NEW style (numFmt 0.0, font9 + color theme1, no border)
$ws.Cells.Item(12, 18).Value = 0.0018411781330637848
$ws.Cells.Item(12, 18).NumberFormat = "0.0"
$ws.Cells.Item(12, 18).Font.Name = "Times New Roman"
$ws.Cells.Item(12, 18).Font.Size = 9
$ws.Cells.Item(12, 18).Font.ThemeColor = 1

# Row 13: "-" text, same style as Q13
$ws.Cells.Item(13, 18).Value = "-"
$ws.Cells.Item(13, 18).Style = $ws.Cells.Item(13, 17).Style

# Update selection to match the authored diff (S4)
$ws.Range("S4").Select()
